$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# New age-sampling counts for B17:B117 (rounding update)
$newValues = @(0, 0, 0, 51644, 57082, 60109, 62268, 62977, 65017, 65748, 67616, 69524, 68853, 70349, 69066, 70911, 70593, 70290, 70329, 68168, 66195, 63651, 61078, 58825, 57758, 57888, 57474, 58945, 58764, 60562, 62188, 58516, 57151, 55721, 52759, 52699, 51203, 52390, 52233, 50728, 50400, 48635, 46799, 46043, 44144, 44053, 41591, 40048, 39347, 38105, 37557, 36069, 35044, 35468, 36916, 30984, 28539, 27857, 24576, 24515, 22998, 22513, 21180, 20299, 19183, 18359, 16639, 15860, 14577, 13563, 13090, 12411, 10625, 9448, 8392, 7237, 5790, 4815, 3856, 2895, 2330, 1584, 903, 684, 480, 260, 184, 107, 62, 40, 29, 9, 5, 6, 0, 3, 0, 0, 4, 3, 0)

$startRow = 17
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws1.Cells.Item($startRow + $i, 2).Value = $newValues[$i]
}

# Rows 20:117 pick up the "format" sheet's General-number style (matches
# the style already used elsewhere in the workbook, e.g. format!A5)
$srcStyle = $ws3.Range("A5")
$srcStyle.Copy()
$ws1.Range("B20:B117").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the selection left behind by the edit
[void]$ws1.Range("B16:B19").Select()
